$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "User1"/"Pass1" row with "desy"/"desy"
$ws.Range("A2").Value = "desy"
$ws.Range("B2").Value = "desy"

# Move the active selection to A3
$ws.Range("A3").Select()
